$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        do {
            $found.Value = "In Translation"
            $found = $used.FindNext($found)
        } while ($found -ne $null -and $found.Address() -ne $firstAddress)
    }
    $ws.Columns.AutoFit() | Out-Null
}
